# Update for latest Morgan poll
# New Morgan poll results go into the "Latest Morgan ->" row (row 6),
# shifting the previous "Latest Morgan" down to "Second Morgan" (row 7)
# and the previous "Second Morgan" down to "Third Morgan" (row 8).
# The previous "Third Morgan" row is discarded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Shift existing poll rows down: row 7 <- old row 6, row 8 <- old row 7
$ws.Range("B8:G8").Value = $ws.Range("B7:G7").Value2
$ws.Range("B7:G7").Value = $ws.Range("B6:G6").Value2

# New latest Morgan poll numbers into row 6
$ws.Range("B6").Value = 55.5
$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 49
$ws.Range("F6").Value = 57
$ws.Range("G6").Value = 63.5

# Update selection to match the saved view state
$ws.Range("G7").Select()
